$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-01-06 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-07 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("28-17=11", $true, $false, $false, $false, $false, $true, 1, $false, "39+55=94", 2) | Out-Null
$d.Content.Find.Execute("95-25=70", $true, $false, $false, $false, $false, $true, 1, $false, "20+30=50", 2) | Out-Null
$d.Content.Find.Execute("14+44=58", $true, $false, $false, $false, $false, $true, 1, $false, "73-6=67", 2) | Out-Null
$d.Content.Find.Execute("99-84=15", $true, $false, $false, $false, $false, $true, 1, $false, "7+90=97", 2) | Out-Null
$d.Content.Find.Execute("32+20=52", $true, $false, $false, $false, $false, $true, 1, $false, "61+33=94", 2) | Out-Null
$d.Content.Find.Execute("55-38=17", $true, $false, $false, $false, $false, $true, 1, $false, "50-16=34", 2) | Out-Null
$d.Content.Find.Execute("73-11=62", $true, $false, $false, $false, $false, $true, 1, $false, "17+82=99", 2) | Out-Null
$d.Content.Find.Execute("4+57=61", $true, $false, $false, $false, $false, $true, 1, $false, "31+53=84", 2) | Out-Null
$d.Content.Find.Execute("70-9=61", $true, $false, $false, $false, $false, $true, 1, $false, "49-7=42", 2) | Out-Null
$d.Content.Find.Execute("17+47=64", $true, $false, $false, $false, $false, $true, 1, $false, "3+55=58", 2) | Out-Null
$d.Content.Find.Execute("35+32=67", $true, $false, $false, $false, $false, $true, 1, $false, "23+53=76", 2) | Out-Null
$d.Content.Find.Execute("82-0=82", $true, $false, $false, $false, $false, $true, 1, $false, "82-28=54", 2) | Out-Null
$d.Content.Find.Execute("23+55=78", $true, $false, $false, $false, $false, $true, 1, $false, "43-18=25", 2) | Out-Null
$d.Content.Find.Execute("25+61=86", $true, $false, $false, $false, $false, $true, 1, $false, "5+70=75", 2) | Out-Null
$d.Content.Find.Execute("88-20=68", $true, $false, $false, $false, $false, $true, 1, $false, "79-61=18", 2) | Out-Null
$d.Content.Find.Execute("37+37=74", $true, $false, $false, $false, $false, $true, 1, $false, "84+0=84", 2) | Out-Null
$d.Content.Find.Execute("48+29=77", $true, $false, $false, $false, $false, $true, 1, $false, "89-25=64", 2) | Out-Null
$d.Content.Find.Execute("47+49=96", $true, $false, $false, $false, $false, $true, 1, $false, "86-37=49", 2) | Out-Null
$d.Content.Find.Execute("78-44=34", $true, $false, $false, $false, $false, $true, 1, $false, "45+8=53", 2) | Out-Null
$d.Content.Find.Execute("35-19=16", $true, $false, $false, $false, $false, $true, 1, $false, "57+8=65", 2) | Out-Null
$d.Content.Find.Execute("41-20=21", $true, $false, $false, $false, $false, $true, 1, $false, "47+30=77", 2) | Out-Null
$d.Content.Find.Execute("50+37=87", $true, $false, $false, $false, $false, $true, 1, $false, "52+24=76", 2) | Out-Null
$d.Content.Find.Execute("72-21=51", $true, $false, $false, $false, $false, $true, 1, $false, "24+37=61", 2) | Out-Null
$d.Content.Find.Execute("46+44=90", $true, $false, $false, $false, $false, $true, 1, $false, "60-1=59", 2) | Out-Null
$d.Content.Find.Execute("89-26=63", $true, $false, $false, $false, $false, $true, 1, $false, "1+56=57", 2) | Out-Null
$d.Content.Find.Execute("69+15=84", $true, $false, $false, $false, $false, $true, 1, $false, "84-46=38", 2) | Out-Null
$d.Content.Find.Execute("8-7=1", $true, $false, $false, $false, $false, $true, 1, $false, "80+17=97", 2) | Out-Null
$d.Content.Find.Execute("46-24=22", $true, $false, $false, $false, $false, $true, 1, $false, "38+36=74", 2) | Out-Null
$d.Content.Find.Execute("82-64=18", $true, $false, $false, $false, $false, $true, 1, $false, "6+8=14", 2) | Out-Null
$d.Content.Find.Execute("57-0=57", $true, $false, $false, $false, $false, $true, 1, $false, "6+82=88", 2) | Out-Null
$d.Content.Find.Execute("88-8=80", $true, $false, $false, $false, $false, $true, 1, $false, "57-43=14", 2) | Out-Null
$d.Content.Find.Execute("2+17=19", $true, $false, $false, $false, $false, $true, 1, $false, "6+86=92", 2) | Out-Null
$d.Content.Find.Execute("60-46=14", $true, $false, $false, $false, $false, $true, 1, $false, "71-20=51", 2) | Out-Null
$d.Content.Find.Execute("18-0=18", $true, $false, $false, $false, $false, $true, 1, $false, "46+28=74", 2) | Out-Null
$d.Content.Find.Execute("34-21=13", $true, $false, $false, $false, $false, $true, 1, $false, "59+18=77", 2) | Out-Null
$d.Content.Find.Execute("37+7=44", $true, $false, $false, $false, $false, $true, 1, $false, "0+76=76", 2) | Out-Null
$d.Content.Find.Execute("4+59=63", $true, $false, $false, $false, $false, $true, 1, $false, "46-8=38", 2) | Out-Null
$d.Content.Find.Execute("18+66=84", $true, $false, $false, $false, $false, $true, 1, $false, "35-15=20", 2) | Out-Null
$d.Content.Find.Execute("83+4=87", $true, $false, $false, $false, $false, $true, 1, $false, "31+53=84", 2) | Out-Null
$d.Content.Find.Execute("98-33=65", $true, $false, $false, $false, $false, $true, 1, $false, "3+73=76", 2) | Out-Null
$d.Content.Find.Execute("80-73=7", $true, $false, $false, $false, $false, $true, 1, $false, "86-17=69", 2) | Out-Null
$d.Content.Find.Execute("61-60=1", $true, $false, $false, $false, $false, $true, 1, $false, "90-32=58", 2) | Out-Null
$d.Content.Find.Execute("77-0=77", $true, $false, $false, $false, $false, $true, 1, $false, "77-61=16", 2) | Out-Null
$d.Content.Find.Execute("40+35=75", $true, $false, $false, $false, $false, $true, 1, $false, "36+59=95", 2) | Out-Null
$d.Content.Find.Execute("15+45=60", $true, $false, $false, $false, $false, $true, 1, $false, "16+25=41", 2) | Out-Null
$d.Content.Find.Execute("16+77=93", $true, $false, $false, $false, $false, $true, 1, $false, "71-53=18", 2) | Out-Null
$d.Content.Find.Execute("79-28=51", $true, $false, $false, $false, $false, $true, 1, $false, "73-69=4", 2) | Out-Null
$d.Content.Find.Execute("75-5=70", $true, $false, $false, $false, $false, $true, 1, $false, "82-78=4", 2) | Out-Null
$d.Content.Find.Execute("24+5=29", $true, $false, $false, $false, $false, $true, 1, $false, "38-19=19", 2) | Out-Null
$d.Content.Find.Execute("42+56=98", $true, $false, $false, $false, $false, $true, 1, $false, "41-26=15", 2) | Out-Null
$d.Content.Find.Execute("83-38=45", $true, $false, $false, $false, $false, $true, 1, $false, "89-7=82", 2) | Out-Null
$d.Content.Find.Execute("57-31=26", $true, $false, $false, $false, $false, $true, 1, $false, "21+10=31", 2) | Out-Null
$d.Content.Find.Execute("87-11=76", $true, $false, $false, $false, $false, $true, 1, $false, "24+47=71", 2) | Out-Null
$d.Content.Find.Execute("13+77=90", $true, $false, $false, $false, $false, $true, 1, $false, "79+8=87", 2) | Out-Null
$d.Content.Find.Execute("40-35=5", $true, $false, $false, $false, $false, $true, 1, $false, "18-8=10", 2) | Out-Null
$d.Content.Find.Execute("92-13=79", $true, $false, $false, $false, $false, $true, 1, $false, "96-65=31", 2) | Out-Null
$d.Content.Find.Execute("86-42=44", $true, $false, $false, $false, $false, $true, 1, $false, "54+45=99", 2) | Out-Null
$d.Content.Find.Execute("95-78=17", $true, $false, $false, $false, $false, $true, 1, $false, "39+17=56", 2) | Out-Null
$d.Content.Find.Execute("83-25=58", $true, $false, $false, $false, $false, $true, 1, $false, "7+33=40", 2) | Out-Null
$d.Content.Find.Execute("14+81=95", $true, $false, $false, $false, $false, $true, 1, $false, "16-1=15", 2) | Out-Null
$d.Content.Find.Execute("61-55=6", $true, $false, $false, $false, $false, $true, 1, $false, "84-27=57", 2) | Out-Null
$d.Content.Find.Execute("95-72=23", $true, $false, $false, $false, $false, $true, 1, $false, "50-23=27", 2) | Out-Null
$d.Content.Find.Execute("32+47=79", $true, $false, $false, $false, $false, $true, 1, $false, "3+50=53", 2) | Out-Null
$d.Content.Find.Execute("13+4=17", $true, $false, $false, $false, $false, $true, 1, $false, "47-6=41", 2) | Out-Null
$d.Content.Find.Execute("53-17=36", $true, $false, $false, $false, $false, $true, 1, $false, "29+60=89", 2) | Out-Null
$d.Content.Find.Execute("54-50=4", $true, $false, $false, $false, $false, $true, 1, $false, "46+21=67", 2) | Out-Null
$d.Content.Find.Execute("51-2=49", $true, $false, $false, $false, $false, $true, 1, $false, "94-47=47", 2) | Out-Null
$d.Content.Find.Execute("82-74=8", $true, $false, $false, $false, $false, $true, 1, $false, "14+48=62", 2) | Out-Null
$d.Content.Find.Execute("35-9=26", $true, $false, $false, $false, $false, $true, 1, $false, "61-35=26", 2) | Out-Null
$d.Content.Find.Execute("0+46=46", $true, $false, $false, $false, $false, $true, 1, $false, "17-15=2", 2) | Out-Null
$d.Content.Find.Execute("11+13=24", $true, $false, $false, $false, $false, $true, 1, $false, "59-44=15", 2) | Out-Null
$d.Content.Find.Execute("16+30=46", $true, $false, $false, $false, $false, $true, 1, $false, "71+19=90", 2) | Out-Null
$d.Content.Find.Execute("77-4=73", $true, $false, $false, $false, $false, $true, 1, $false, "20+34=54", 2) | Out-Null
$d.Content.Find.Execute("11+42=53", $true, $false, $false, $false, $false, $true, 1, $false, "52+35=87", 2) | Out-Null
$d.Content.Find.Execute("19+11=30", $true, $false, $false, $false, $false, $true, 1, $false, "77+7=84", 2) | Out-Null
$d.Content.Find.Execute("66-52=14", $true, $false, $false, $false, $false, $true, 1, $false, "14+41=55", 2) | Out-Null
$d.Content.Find.Execute("47+27=74", $true, $false, $false, $false, $false, $true, 1, $false, "90-22=68", 2) | Out-Null
$d.Content.Find.Execute("88-59=29", $true, $false, $false, $false, $false, $true, 1, $false, "90-62=28", 2) | Out-Null
$d.Content.Find.Execute("43+8=51", $true, $false, $false, $false, $false, $true, 1, $false, "54-2=52", 2) | Out-Null
$d.Content.Find.Execute("31+3=34", $true, $false, $false, $false, $false, $true, 1, $false, "50+10=60", 2) | Out-Null
$d.Content.Find.Execute("74+17=91", $true, $false, $false, $false, $false, $true, 1, $false, "76+15=91", 2) | Out-Null
$d.Content.Find.Execute("17+4=21", $true, $false, $false, $false, $false, $true, 1, $false, "37+41=78", 2) | Out-Null
$d.Content.Find.Execute("56-23=33", $true, $false, $false, $false, $false, $true, 1, $false, "40+38=78", 2) | Out-Null
$d.Content.Find.Execute("70+5=75", $true, $false, $false, $false, $false, $true, 1, $false, "7+39=46", 2) | Out-Null
$d.Content.Find.Execute("11+41=52", $true, $false, $false, $false, $false, $true, 1, $false, "23+42=65", 2) | Out-Null
$d.Content.Find.Execute("59-6=53", $true, $false, $false, $false, $false, $true, 1, $false, "50-5=45", 2) | Out-Null
$d.Content.Find.Execute("79+2=81", $true, $false, $false, $false, $false, $true, 1, $false, "76+1=77", 2) | Out-Null
$d.Content.Find.Execute("52+4=56", $true, $false, $false, $false, $false, $true, 1, $false, "81-19=62", 2) | Out-Null
$d.Content.Find.Execute("49-2=47", $true, $false, $false, $false, $false, $true, 1, $false, "12+83=95", 2) | Out-Null
$d.Content.Find.Execute("70-69=1", $true, $false, $false, $false, $false, $true, 1, $false, "74-56=18", 2) | Out-Null
$d.Content.Find.Execute("91-11=80", $true, $false, $false, $false, $false, $true, 1, $false, "88-4=84", 2) | Out-Null
$d.Content.Find.Execute("51+4=55", $true, $false, $false, $false, $false, $true, 1, $false, "34+62=96", 2) | Out-Null
$d.Content.Find.Execute("18+51=69", $true, $false, $false, $false, $false, $true, 1, $false, "88-87=1", 2) | Out-Null
$d.Content.Find.Execute("9+30=39", $true, $false, $false, $false, $false, $true, 1, $false, "15-12=3", 2) | Out-Null
$d.Content.Find.Execute("22+68=90", $true, $false, $false, $false, $false, $true, 1, $false, "78-46=32", 2) | Out-Null
$d.Content.Find.Execute("27+3=30", $true, $false, $false, $false, $false, $true, 1, $false, "43-39=4", 2) | Out-Null
$d.Content.Find.Execute("50-30=20", $true, $false, $false, $false, $false, $true, 1, $false, "81-33=48", 2) | Out-Null
$d.Content.Find.Execute("50+32=82", $true, $false, $false, $false, $false, $true, 1, $false, "8+68=76", 2) | Out-Null
$d.Content.Find.Execute("81-22=59", $true, $false, $false, $false, $false, $true, 1, $false, "58-16=42", 2) | Out-Null
$d.Content.Find.Execute("65+31=96", $true, $false, $false, $false, $false, $true, 1, $false, "24+53=77", 2) | Out-Null

Write-Output "Replacements complete"
